# The only content-visible change in this revision is that the paragraph
# holding the first row of the ITRF covariance matrix (the paragraph that
# immediately follows the "Covariance in ITRF" line) comes back from the
# spell-checker clean, so Word stamps it w14:noSpellErr="1" on the <w:p>.
#
# Locate that paragraph robustly (by text, not a hard-coded index), then
# round-trip its WordOpenXML through a minimal regex patch that adds the
# w14:noSpellErr="1" attribute to its <w:p> element, and write it back with
# InsertXML (InsertXML replaces the contents of the exact range it's called
# on, so calling it on just this paragraph's Range leaves everything else
# untouched).

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "in ITRF") {
        $target = $p.Next()
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $xml = $r.WordOpenXML
    $newXml = [regex]::Replace($xml, '<w:p([ >])', '<w:p w14:noSpellErr="1"$1', 1)
    $r.InsertXML($newXml)
}
